# Add two new rows documenting the new "crop_segmentation" and "ave_ref_crop"
# functions just above the existing "transformation.py" entry, shifting the
# rows below down by two (matching the target diff).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at row 16 (pushes "transformation.py" and everything
# after it down to rows 18+).
$null = $ws.Rows("16:17").Insert()

# Row 16: crop_segmentation / f
$ws.Range("E16").Value2 = "crop_segmentation"
$ws.Range("F16").Value2 = "f"

# Row 17: ave_ref_crop / f, average the reflectance of crops
$ws.Range("E17").Value2 = "ave_ref_crop"
$ws.Range("F17").Value2 = "f, average the reflectance of crops"

# Update the sheet's active selection to F17, matching the saved workbook view.
$null = $ws.Range("F17").Select()
